# Add WPF tutorial progress row to the Pluralsight tracker sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New course entry (row 3): course title + its current status.
$ws.Range("A3").Value = "Building an Enterprise App with WPF, MVVM, and Entity Framework Code First"
$ws.Range("B3").Value = "CH 4 - START"

# Column A now holds a much longer course title - widen it to fit.
$ws.Columns.Item(1).AutoFit() | Out-Null

# Leave the selection where Excel would land after typing into B3 and hitting Enter.
$ws.Range("B4").Select() | Out-Null
